# Add missing boards (15, 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Board 15 (row 16) - Chip ID/Flash ID/Mac Address were missing
$ws.Range("B16").Value = "Lite"
$ws.Range("C16").Value = "0x000a7262"
$ws.Range("D16").Value = "0x4016"
$ws.Range("E16").Value = "5C:CF:7F:A:72:62"

# Board 9 (row 10) - Chip ID/Flash ID/Mac Address were missing
$ws.Range("B10").Value = "Lite"
$ws.Range("C10").Value = "0x00d76690"
$ws.Range("D10").Value = "0x4016"
$ws.Range("E10").Value = "18:FE:34:D7:66:90"

# Leave the active cell selection where the user ended up editing
$ws.Range("C11").Select()
